$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new meeting-diary row for the Oct 26 night meeting -----------
# Columns: A=Date, B=Time start, C=Time end, D=Members present, E=Discussions
$ws.Cells.Item(24,1).Value = 45225                        # 10/26/2023
$ws.Cells.Item(24,2).Value = 0.88888888888888884          # 21:20
$ws.Cells.Item(24,3).Value = 0.93055555555555547          # 22:20
$ws.Cells.Item(24,4).Value = $ws.Cells.Item(23,4).Value2  # same members as row 23
$ws.Cells.Item(24,5).Value = "- Practice presentation part within our team and other team`n- Record the presentation with other teams "

# Copy the formatting "shape" of the row above (date/time/members/discussion
# styles) down onto the freshly-populated row.
$ws.Range("A23:E23").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the row height Excel would have auto-fit for the wrapped text.
$ws.Rows.Item(24).RowHeight = 51

# --- Update the view: scroll down and move the active selection -----------
$ws.Range("A19").Select()
$ws.Range("B25").Select()
